# Apply "bar chart updated (grouped)" edit
# 1) Resumen!C2 and Metricas!B2 : max time value updated
# 2) Metricas!B3 : Z2 time value updated
# 3) Solucion!A2:B41 : Pedido/Salida assignment re-shuffled

$wb = $excel.ActiveWorkbook

$wsResumen = $wb.Worksheets.Item("Resumen")
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsMetricas = $wb.Worksheets.Item("Metricas")

$wsResumen.Range("C2").Value = 640.4113147367285

$wsMetricas.Range("B2").Value = 640.4113147367285
$wsMetricas.Range("B3").Value = 493.0490863877178

$solucionData = @(
  @("Pedido_20", "S001"),
  @("Pedido_11", "S025"),
  @("Pedido_40", "S005"),
  @("Pedido_7", "S029"),
  @("Pedido_19", "S002"),
  @("Pedido_31", "S026"),
  @("Pedido_36", "S006"),
  @("Pedido_1", "S030"),
  @("Pedido_6", "S003"),
  @("Pedido_16", "S027"),
  @("Pedido_26", "S007"),
  @("Pedido_14", "S031"),
  @("Pedido_24", "S004"),
  @("Pedido_30", "S008"),
  @("Pedido_13", "S028"),
  @("Pedido_27", "S009"),
  @("Pedido_35", "S032"),
  @("Pedido_12", "S013"),
  @("Pedido_8", "S010"),
  @("Pedido_17", "S033"),
  @("Pedido_3", "S014"),
  @("Pedido_15", "S037"),
  @("Pedido_21", "S011"),
  @("Pedido_10", "S015"),
  @("Pedido_9", "S034"),
  @("Pedido_4", "S012"),
  @("Pedido_28", "S038"),
  @("Pedido_33", "S016"),
  @("Pedido_25", "S035"),
  @("Pedido_32", "S017"),
  @("Pedido_18", "S039"),
  @("Pedido_29", "S021"),
  @("Pedido_23", "S018"),
  @("Pedido_37", "S036"),
  @("Pedido_2", "S040"),
  @("Pedido_22", "S022"),
  @("Pedido_39", "S019"),
  @("Pedido_34", "S023"),
  @("Pedido_5", "S020"),
  @("Pedido_38", "S024")
)

for ($i = 0; $i -lt $solucionData.Length; $i++) {
  $targetRow = $i + 2
  $pair = $solucionData[$i]
  $wsSolucion.Cells.Item($targetRow, 1).Value = $pair[0]
  $wsSolucion.Cells.Item($targetRow, 2).Value = $pair[1]
}
